# fix: tamplate laporan penjualan
# - Relabel the "total" header to "Total"
# - Center-align the underlying bordered header style (row 5)
# - Widen column H and give column J an explicit width
# - Move the active selection to J12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize the "total" header in I5 -> "Total"
$ws.Range("I5").Value = "Total"

# The bordered header row (A5:L5) gets centered text as well as its border
$ws.Range("A5:L5").HorizontalAlignment = -4108   # xlCenter

# Column width tweaks: H loses its auto bestFit and becomes a touch wider,
# and J gets an explicit width to match the other wide columns
$ws.Columns("H").ColumnWidth = 6.92
$ws.Columns("J").ColumnWidth = 12.42

# Move the selection cursor to J12
$ws.Range("J12").Select()
